$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H27").Value = 10000
$ws.Range("I27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 30000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = -29899
$ws.Range("N27").Value = -30202
$ws.Range("H38").Value = 660.2143
$ws.Range("I38").Value = 155.375
$ws.Range("J38").Value = 1333.3334
$ws.Range("K38").Value = 466.125
$ws.Range("L38").Value = 4000.0002
$ws.Range("M38").Value = -94.125
$ws.Range("N38").Value = -4744.0002
$ws.Range("H62").Value = 2999
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2999
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 2999
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14995
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H106").Value = 38481096
$ws.Range("I106").Value = 41684480
$ws.Range("K106").Value = 41684480
$ws.Range("M106").Value = -41683849
$ws.Range("H112").Value = 1770.4706
$ws.Range("J112").Value = 1932
$ws.Range("L112").Value = 5796
$ws.Range("N112").Value = -8012
$ws.Range("H132").Value = 2434.3845
$ws.Range("I132").Value = 1408.4783
$ws.Range("J132").Value = 10299.667
$ws.Range("K132").Value = 4225.4349
$ws.Range("L132").Value = 30899.001
$ws.Range("M132").Value = -1695.4349
$ws.Range("N132").Value = -35959.001
$ws.Range("H135").Value = 2279.4167
$ws.Range("I135").Value = 1960.875
$ws.Range("K135").Value = 17647.875
$ws.Range("M135").Value = -15112.875
$ws.Range("H138").Value = 2152.6553
$ws.Range("I138").Value = 1907
$ws.Range("J138").Value = 2238.3489
$ws.Range("K138").Value = 5721
$ws.Range("L138").Value = 6715.0467
$ws.Range("M138").Value = -581
$ws.Range("N138").Value = -16995.0467

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 22324.5
$ws.Range("I55").Value = 4000
$ws.Range("K55").Value = 4000
$ws.Range("M55").Value = -3685
$ws.Range("H61").Value = 3908.5789
$ws.Range("I61").Value = 4074.2942
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 4074.2942
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -3862.2942
$ws.Range("N61").Value = -2924
$ws.Range("H110").Value = 4833508.5
$ws.Range("I110").Value = 7409326.5
$ws.Range("J110").Value = 3850
$ws.Range("K110").Value = 7409326.5
$ws.Range("L110").Value = 3850
$ws.Range("M110").Value = -7407281.5
$ws.Range("N110").Value = -7940
$ws.Range("H132").Value = 4761.3076
$ws.Range("J132").Value = 3490.9092
$ws.Range("L132").Value = 10472.7276
$ws.Range("N132").Value = -15532.7276
$ws.Range("H136").Value = 3908.5789
$ws.Range("I136").Value = 4074.2942
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 12222.8826
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -9672.882599999999
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3889.5881
$ws.Range("I20").Value = 2982.818
$ws.Range("K20").Value = 2982.818
$ws.Range("M20").Value = -2735.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 50003620
$ws.Range("J16").Value = 3679.8
$ws.Range("L16").Value = 3679.8
$ws.Range("N16").Value = -4253.8
$ws.Range("H17").Value = 1116.6666
$ws.Range("I17").Value = 950
$ws.Range("J17").Value = 1450
$ws.Range("K17").Value = 950
$ws.Range("L17").Value = 1450
$ws.Range("M17").Value = -776
$ws.Range("N17").Value = -1798
$ws.Range("H22").Value = 7904.4614
$ws.Range("I22").Value = 205.18182
$ws.Range("K22").Value = 205.18182
$ws.Range("M22").Value = 144.81818
$ws.Range("H25").Value = 3979.2856
$ws.Range("I25").Value = 2975.8333
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 2975.8333
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = -2801.8333
$ws.Range("N25").Value = -10348
$ws.Range("H99").Value = 14925.863
$ws.Range("J99").Value = 19062.5
$ws.Range("L99").Value = 19062.5
$ws.Range("N99").Value = -22058.5
$ws.Range("H113").Value = 50003620
$ws.Range("J113").Value = 3679.8
$ws.Range("L113").Value = 3679.8
$ws.Range("N113").Value = -8019.8
$ws.Range("H126").Value = 14925.863
$ws.Range("J126").Value = 19062.5
$ws.Range("L126").Value = 57187.5
$ws.Range("N126").Value = -62127.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1731.6666
$ws.Range("I36").Value = 1731.6666
$ws.Range("K36").Value = 5194.9998
$ws.Range("M36").Value = -5025.9998
$ws.Range("H44").Value = 1428.5714
$ws.Range("J44").Value = 1860
$ws.Range("L44").Value = 5580
$ws.Range("N44").Value = -6376
$ws.Range("H70").Value = 800
$ws.Range("I70").Value = 800
$ws.Range("K70").Value = 2400
$ws.Range("M70").Value = -2085
$ws.Range("H73").Value = 800
$ws.Range("I73").Value = 800
$ws.Range("K73").Value = 2400
$ws.Range("M73").Value = -1308
$ws.Range("H124").Value = 1115
$ws.Range("I124").Value = 1115
$ws.Range("K124").Value = 3345
$ws.Range("M124").Value = 1565
$ws.Range("H132").Value = 2125.6296
$ws.Range("J132").Value = 4400
$ws.Range("L132").Value = 39600
$ws.Range("N132").Value = -44660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 322.6
$ws.Range("J2").Value = 645.6667
$ws.Range("L2").Value = 645.6667
$ws.Range("N2").Value = -871.6667
$ws.Range("H80").Value = 3975
$ws.Range("I80").Value = 1650
$ws.Range("J80").Value = 4750
$ws.Range("K80").Value = 1650
$ws.Range("L80").Value = 4750
$ws.Range("M80").Value = -652
$ws.Range("N80").Value = -6746
$ws.Range("H83").Value = 3975
$ws.Range("I83").Value = 1650
$ws.Range("J83").Value = 4750
$ws.Range("K83").Value = 8250
$ws.Range("L83").Value = 23750
$ws.Range("M83").Value = -3258
$ws.Range("N83").Value = -33734
$ws.Range("H97").Value = 2295.682
$ws.Range("I97").Value = 1678.2142
$ws.Range("K97").Value = 1678.2142
$ws.Range("M97").Value = -1182.2142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 10061.5
$ws.Range("I34").Value = 10000
$ws.Range("J34").Value = 10123
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 10123
$ws.Range("M34").Value = -9828
$ws.Range("N34").Value = -10467
$ws.Range("H46").Value = 3230.5625
$ws.Range("J46").Value = 3307.4167
$ws.Range("L46").Value = 3307.4167
$ws.Range("N46").Value = -3683.4167
$ws.Range("H61").Value = 12347527
$ws.Range("I61").Value = 15874043
$ws.Range("K61").Value = 15874043
$ws.Range("M61").Value = -15873841
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0
$ws.Range("H113").Value = 12347527
$ws.Range("I113").Value = 15874043
$ws.Range("K113").Value = 15874043
$ws.Range("M113").Value = -15871873

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 2000000
$ws.Range("J26").Value = 2000000
$ws.Range("L26").Value = 2000000
$ws.Range("N26").Value = -2000586
$ws.Range("H132").Value = 3459.4666
$ws.Range("J132").Value = 5400
$ws.Range("L132").Value = 16200
$ws.Range("N132").Value = -21260
$ws.Range("H136").Value = 7366.75
$ws.Range("I136").Value = 14998
$ws.Range("K136").Value = 44994
$ws.Range("M136").Value = -42444
